$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert new blank rows at the three locations where new BOM lines were
#    added. Done from the bottom of the sheet upward so that the row
#    numbers used for each Insert() call stay valid (rows below an Insert
#    point are pushed down, rows above are unaffected).
# ---------------------------------------------------------------------------

# New rows for "Button" and "SPDT Switch" (were inserted right before the
# old USB-connector row, i.e. original row 28).
$ws.Range("28:29").Insert()

# New rows for the three TCXO "Clock" options (inserted right before the
# old blank separator that used to sit at row 22, just after Crystal 32kHz).
$ws.Range("22:24").Insert()

# New rows for Resistor 33R/127R/82R5, two separator rows, and the RF
# inductor / ferrite bead pair (inserted right before the old Crystal 8MHz
# row, which used to be row 20).
$ws.Range("20:27").Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the brand-new rows with their BOM data.
# ---------------------------------------------------------------------------

# Row 19 (was an empty separator row, now the "Resistor 390R" line)
$ws.Range("A19").Value = "Resistor 390R"
$ws.Range("B19").Value = "Thin film resistor"
$ws.Range("C19").Value = "594-MCT06030C3900FP5"
$ws.Range("D19").Value = 0.063
$ws.Range("E19").Value = 20
$ws.Range("F19").Formula = "=D19*E19"

# Row 20 - Resistor 33R
$ws.Range("A20").Value = "Resistor 33R"
$ws.Range("B20").Value = "Thin film resistor"
$ws.Range("C20").Value = "603-RT0603FRE0733RL"
$ws.Range("D20").Value = 0.029
$ws.Range("E20").Value = 10
$ws.Range("F20").Formula = "=D20*E20"

# Row 21 - Resistor 127R
$ws.Range("A21").Value = "Resistor 127R"
$ws.Range("B21").Value = "Thin film resistor"
$ws.Range("C21").Value = "603-RT0603FRE07127RL"
$ws.Range("D21").Value = 0.05
$ws.Range("E21").Value = 10
$ws.Range("F21").Formula = "=D21*E21"

# Row 22 - Resistor 82R5
$ws.Range("A22").Value = "Resistor 82R5"
$ws.Range("B22").Value = "Thin film resistor"
$ws.Range("C22").Value = "594-MCT06030C8259FP5"
$ws.Range("D22").Value = 0.115
$ws.Range("E22").Value = 10
$ws.Range("F22").Formula = "=D22*E22"

# Row 23 - empty separator row, but keeps a formula (D/E blank -> 0)
$ws.Range("F23").Formula = "=D23*E23"

# Row 24 - empty separator row (no formula at all)

# Row 25 - RF inductor
$ws.Range("A25").Value = "RF inductor"
$ws.Range("C25").Value = "81-LQG15HH7N5G02D"
$ws.Range("D25").Value = 0.184
$ws.Range("E25").Value = 10
$ws.Range("F25").Formula = "=D25*E25"

# Row 26 - Ferrite bead
$ws.Range("A26").Value = "Ferrite bead"
$ws.Range("C26").Value = "81-BLM18KG121TN1D"
$ws.Range("D26").Value = 0.057
$ws.Range("E26").Value = 10
$ws.Range("F26").Formula = "=D26*E26"

# Row 27 - empty separator row (no formula)

# Row 28 - Crystal 8MHz entry: footprint / part number change
$ws.Range("B28").Value = "Crystals 8.000 MHz 8 pF"
$ws.Range("C28").Value = "344-NX3225GD8MCRA3"
$ws.Range("D28").Value = 0.645

# Row 29 (Crystal 32kHz) unchanged

# Row 30 - Clock option 1 (LVPECL)
$ws.Range("A30").Value = "Clock"
$ws.Range("C30").Value = "815-AX3PAF1-125.0000"
$ws.Range("D30").Value = 5.14
$ws.Range("E30").Value = 1
$ws.Range("F30").Formula = "=D30*E30"

# Row 31 - Clock option 2 (HCSL)
$ws.Range("A31").Value = "Clock"
$ws.Range("C31").Value = "815-AX3HAF1135.0000T"
$ws.Range("D31").Value = 4.61
$ws.Range("E31").Value = 1
$ws.Range("F31").Formula = "=D31*E31"

# Row 32 - Clock option 3 (LVDS)
$ws.Range("A32").Value = "Clock"
$ws.Range("C32").Value = "815-AX3DAF1-125.0000"
$ws.Range("D32").Value = 5.42
$ws.Range("E32").Value = 1
$ws.Range("F32").Formula = "=D32*E32"

# Row 39 - Button
$ws.Range("A39").Value = "Button"
$ws.Range("C39").Value = "710-434781025816"
$ws.Range("D39").Value = 1.08
$ws.Range("E39").Value = 3
$ws.Range("F39").Formula = "=D39*E39"

# Row 40 - SPDT Switch
$ws.Range("A40").Value = "SPDT Switch"
$ws.Range("C40").Value = "611-7101TPCBE"
$ws.Range("D40").Value = 15.35
$ws.Range("E40").Value = 1
$ws.Range("F40").Formula = "=D40*E40"
